$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# ALC row 74
$ws_ALC.Range("H74").Value = 3816.327
$ws_ALC.Range("I74").Value = 3809.9714
$ws_ALC.Range("J74").Value = 3829.4119
$ws_ALC.Range("K74").Value = 3809.9714
$ws_ALC.Range("L74").Value = 3829.4119
$ws_ALC.Range("M74").Value = -2873.9714
$ws_ALC.Range("N74").Value = -5701.4119

# ALC row 75
$ws_ALC.Range("H75").Value = 43000
$ws_ALC.Range("I75").Value = 0
$ws_ALC.Range("J75").Value = 43000
$ws_ALC.Range("K75").Value = 0
$ws_ALC.Range("L75").Value = 43000
$ws_ALC.Range("M75").ClearContents()
$ws_ALC.Range("N75").Value = -44872

# ALC row 76
$ws_ALC.Range("H76").Value = 3068.1333
$ws_ALC.Range("I76").Value = 2994.889
$ws_ALC.Range("J76").Value = 3178
$ws_ALC.Range("K76").Value = 2994.889
$ws_ALC.Range("L76").Value = 3178
$ws_ALC.Range("M76").Value = -2679.889
$ws_ALC.Range("N76").Value = -3808

# ALC row 77
$ws_ALC.Range("H77").Value = 3816.327
$ws_ALC.Range("I77").Value = 3809.9714
$ws_ALC.Range("J77").Value = 3829.4119
$ws_ALC.Range("K77").Value = 19049.857
$ws_ALC.Range("L77").Value = 19147.0595
$ws_ALC.Range("M77").Value = -14369.857
$ws_ALC.Range("N77").Value = -28507.0595

# ALC row 78
$ws_ALC.Range("H78").Value = 43000
$ws_ALC.Range("I78").Value = 0
$ws_ALC.Range("J78").Value = 43000
$ws_ALC.Range("K78").Value = 0
$ws_ALC.Range("L78").Value = 129000
$ws_ALC.Range("M78").ClearContents()
$ws_ALC.Range("N78").Value = -138360

# ALC row 79
$ws_ALC.Range("H79").Value = 3068.1333
$ws_ALC.Range("I79").Value = 2994.889
$ws_ALC.Range("J79").Value = 3178
$ws_ALC.Range("K79").Value = 2994.889
$ws_ALC.Range("L79").Value = 3178
$ws_ALC.Range("M79").Value = -1902.889
$ws_ALC.Range("N79").Value = -5362

# ALC row 80
$ws_ALC.Range("H80").Value = 10732855
$ws_ALC.Range("I80").Value = 18519122
$ws_ALC.Range("J80").Value = 721939.2
$ws_ALC.Range("K80").Value = 55557366
$ws_ALC.Range("L80").Value = 2165817.6
$ws_ALC.Range("M80").Value = -55556368
$ws_ALC.Range("N80").Value = -2167813.6

# ALC row 83
$ws_ALC.Range("H83").Value = 10732855
$ws_ALC.Range("I83").Value = 18519122
$ws_ALC.Range("J83").Value = 721939.2
$ws_ALC.Range("K83").Value = 166672098
$ws_ALC.Range("L83").Value = 6497452.8
$ws_ALC.Range("M83").Value = -166667106
$ws_ALC.Range("N83").Value = -6507436.8

# ALC row 86
$ws_ALC.Range("H86").Value = 6254071
$ws_ALC.Range("I86").Value = 12503792
$ws_ALC.Range("K86").Value = 12503792
$ws_ALC.Range("M86").Value = -12502669

# ALC row 87
$ws_ALC.Range("H87").Value = 19840.5
$ws_ALC.Range("J87").Value = 19840.5
$ws_ALC.Range("L87").Value = 19840.5
$ws_ALC.Range("N87").Value = -22336.5

# ALC row 89
$ws_ALC.Range("H89").Value = 6254071
$ws_ALC.Range("I89").Value = 12503792
$ws_ALC.Range("K89").Value = 62518960
$ws_ALC.Range("M89").Value = -62513344

# ALC row 90
$ws_ALC.Range("H90").Value = 19840.5
$ws_ALC.Range("J90").Value = 19840.5
$ws_ALC.Range("L90").Value = 59521.5
$ws_ALC.Range("N90").Value = -72001.5

# ALC row 93
$ws_ALC.Range("H93").Value = 0
$ws_ALC.Range("J93").Value = 0
$ws_ALC.Range("L93").Value = 0
$ws_ALC.Range("N93").ClearContents()

# ALC row 96
$ws_ALC.Range("H96").Value = 84814.57000000001
$ws_ALC.Range("I96").Value = 4972.875
$ws_ALC.Range("J96").Value = 191270.17
$ws_ALC.Range("K96").Value = 14918.625
$ws_ALC.Range("L96").Value = 573810.51
$ws_ALC.Range("M96").Value = -13545.625
$ws_ALC.Range("N96").Value = -576556.51

# ALC row 107
$ws_ALC.Range("H107").Value = 7610.2856
$ws_ALC.Range("I107").Value = 8545.333000000001
$ws_ALC.Range("J107").Value = 2000
$ws_ALC.Range("K107").Value = 8545.333000000001
$ws_ALC.Range("L107").Value = 2000
$ws_ALC.Range("M107").Value = -6625.333000000001
$ws_ALC.Range("N107").Value = -5840

# ALC row 133
$ws_ALC.Range("H133").Value = 56515.8
$ws_ALC.Range("J133").Value = 56515.8
$ws_ALC.Range("L133").Value = 56515.8
$ws_ALC.Range("N133").Value = -66635.8

# ALC row 137
$ws_ALC.Range("H137").Value = 3246.204
$ws_ALC.Range("I137").Value = 864.34375
$ws_ALC.Range("J137").Value = 7729.706
$ws_ALC.Range("K137").Value = 2593.03125
$ws_ALC.Range("L137").Value = 23189.118
$ws_ALC.Range("M137").Value = -43.03125
$ws_ALC.Range("N137").Value = -28289.118

# ALC row 138
$ws_ALC.Range("H138").Value = 1590.1414
$ws_ALC.Range("I138").Value = 809.3125
$ws_ALC.Range("J138").Value = 1963.0746
$ws_ALC.Range("K138").Value = 2427.9375
$ws_ALC.Range("L138").Value = 5889.2238
$ws_ALC.Range("M138").Value = 2712.0625
$ws_ALC.Range("N138").Value = -16169.2238

# ALC row 139
$ws_ALC.Range("H139").Value = 42275
$ws_ALC.Range("J139").Value = 42275
$ws_ALC.Range("L139").Value = 42275
$ws_ALC.Range("N139").Value = -52555

# ALC row 140
$ws_ALC.Range("H140").Value = 49499.5
$ws_ALC.Range("J140").Value = 49499.5
$ws_ALC.Range("L140").Value = 49499.5
$ws_ALC.Range("N140").Value = -59859.5

# ALC row 141
$ws_ALC.Range("H141").Value = 3326.3333
$ws_ALC.Range("I141").Value = 1244.9524
$ws_ALC.Range("J141").Value = 6968.75
$ws_ALC.Range("K141").Value = 3734.857199999999
$ws_ALC.Range("L141").Value = 20906.25
$ws_ALC.Range("M141").Value = 1445.142800000001
$ws_ALC.Range("N141").Value = -31266.25

# CUL row 5
$ws_CUL.Range("H5").Value = 4809.32
$ws_CUL.Range("J5").Value = 1640.7142
$ws_CUL.Range("L5").Value = 4922.142599999999
$ws_CUL.Range("N5").Value = -5146.142599999999

# CUL row 103
$ws_CUL.Range("H103").Value = 7575
$ws_CUL.Range("J103").Value = 0
$ws_CUL.Range("L103").Value = 0
$ws_CUL.Range("N103").ClearContents()

# CUL row 106
$ws_CUL.Range("H106").Value = 260000
$ws_CUL.Range("J106").Value = 0
$ws_CUL.Range("L106").Value = 0
$ws_CUL.Range("N106").ClearContents()

# CUL row 132
$ws_CUL.Range("H132").Value = 2581.8333
$ws_CUL.Range("I132").Value = 2019.8
$ws_CUL.Range("J132").Value = 2798
$ws_CUL.Range("K132").Value = 18178.2
$ws_CUL.Range("L132").Value = 25182
$ws_CUL.Range("M132").Value = -15648.2
$ws_CUL.Range("N132").Value = -30242

# CUL row 135
$ws_CUL.Range("H135").Value = 4809.32
$ws_CUL.Range("J135").Value = 1640.7142
$ws_CUL.Range("L135").Value = 14766.4278
$ws_CUL.Range("N135").Value = -19836.4278

# GSM row 113
$ws_GSM.Range("H113").Value = 1799.8182
$ws_GSM.Range("I113").Value = 1699.7778
$ws_GSM.Range("J113").Value = 2250
$ws_GSM.Range("K113").Value = 1699.7778
$ws_GSM.Range("L113").Value = 2250
$ws_GSM.Range("M113").Value = 470.2221999999999
$ws_GSM.Range("N113").Value = -6590

# LTW row 122
$ws_LTW.Range("H122").Value = 145013.28
$ws_LTW.Range("I122").Value = 335034.66
$ws_LTW.Range("J122").Value = 2497.25
$ws_LTW.Range("K122").Value = 1005103.98
$ws_LTW.Range("L122").Value = 7491.75
$ws_LTW.Range("M122").Value = -1002653.98
$ws_LTW.Range("N122").Value = -12391.75
